$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they keep the same bold/border/alignment style.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 94
    $ws.Cells.Item($r, 31).Value = 68
    $ws.Cells.Item($r, 32).Value = 0
}
